$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab to reflect the new "through" date
$ws.Name = "Through 2022-03-29"

# Update the header label in I1 (shared string "2022 (through 03-28)" -> "2022 (through 03-29)")
$ws.Range("I1").Value = "2022 (through 03-29)"

# Update the March (row 4) and Total (row 14) figures for the "2022" column (I)
$ws.Range("I4").Value = 125
$ws.Range("I14").Value = 425
